$wb = $excel.ActiveWorkbook

# Rename existing third sheet (2020_02_23 -> 2020_02_24) and append a new
# sheet for 2020_02_25 right after it, mirroring the day-by-day layout of
# the workbook.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "2020_02_24_covid19_country"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "2020_02_25_covid19_country"

$dataRows = @(
  ,('country', 'total cases', 'new cases', 'total deaths', 'total case', 'Total recover', 'Critical', 'region')
  ,('China', 77666, 516, 2664, 72, 27636, 9126, 'Asia ')
  ,('South Korea', 977, 144, 11, 3, 22, 6, 'Asia')
  ,('Diamon Princess', 691, 0, 4, 1, 10, 35, 'Asia (at the Japanese harbour) ')
  ,('Italy', 287, 58, 7, 0, 1, 23, 'Europe')
  ,('Japan', 161, 2, 1, 0, 23, 7, 'Asia')
  ,('Iran', 95, 34, 15, 3, 3, 0, 'Middle East Asia')
  ,('Singapore', 90, 0, 0, 0, 53, 7, 'South East Asia')
  ,('Hong Kong', 81, 0, 2, 0, 12, 6, 'Asia')
  ,('United States', 35, 0, 0, 0, 6, 0, 'North America ')
  ,('Thailand', 37, 2, 0, 0, 21, 2, 'South East Asia')
  ,('Taiwan', 31, 1, 1, 0, 5, 1, 'Asia')
  ,('Australia', 22, 0, 0, 0, 11, 0, 'Australia/Oceania')
  ,('Malaysia', 22, 0, 0, 0, 20, 0, 'Asia')
  ,('Germany', 16, 0, 0, 0, 14, 0, 'Europe')
  ,('Vietnam', 16, 0, 0, 0, 16, 0, 'South East Asia')
  ,('United Kingdom', 13, 0, 0, 0, 8, 0, 'Europe')
  ,('United Arab Emirates', 13, 0, 0, 0, 3, 2, 'Middle East')
  ,('France', 12, 0, 1, 0, 11, 0, 'Europe')
  ,('Canada', 11, 0, 0, 0, 3, 0, 'North America')
  ,('Macao', 10, 0, 0, 0, 5, 0, 'Asia ')
  ,('Bahrain', 8, 6, 0, 0, 0, 0, 'Asaia')
  ,('Kuwait', 8, 3, 0, 0, 0, 0, 'Middle East ')
  ,('Iraq', 5, 4, 0, 0, 0, 0, 'Middle East ')
  ,('Spain', 4, 1, 0, 0, 2, 0, 'Europe')
  ,('Philippines', 3, 0, 1, 0, 2, 0, 'South East Asia')
  ,('India', 3, 0, 0, 0, 3, 0, 'Asia ')
  ,('Austria', 2, 2, 0, 0, 0, 0, 'Europe')
  ,('Israel', 2, 0, 0, 0, 0, 0, 'Middle East')
  ,('Oman', 2, 0, 0, 0, 0, 0, 'Middle East ')
  ,('Russia', 2, 0, 0, 0, 2, 0, 'Europe')
  ,('Afghanistan', 1, 1, 0, 0, 0, 0, 'Asia')
  ,('Cambodia', 1, 0, 0, 0, 1, 0, 'South East Asia')
  ,('Belgium', 1, 0, 0, 0, 1, 0, 'Europe')
  ,('Coratia', 1, 1, 0, 0, 0, 0, 'Europe')
  ,('Finland', 1, 0, 0, 0, 1, 0, 'Europe')
  ,('Egypt', 1, 0, 0, 0, 1, 0, 'Africa ')
  ,('Nepal', 1, 0, 0, 0, 1, 0, 'Asia')
  ,('Lebanon', 1, 0, 0, 0, 0, 0, 'Middle East ')
  ,('Sri Lanka', 1, 0, 0, 0, 1, 0, 'Asia ')
  ,('Sweden', 1, 0, 0, 0, 0, 0, 'Europe')
  ,('Switzerlan', 1, 0, 0, 0, 0, 0, 'Europe')
)

$nRows = $dataRows.Count
$nCols = $dataRows[0].Count
$arr = New-Object 'object[,]' $nRows, $nCols
for ($r = 0; $r -lt $nRows; $r++) {
  for ($c = 0; $c -lt $nCols; $c++) {
    $arr[$r, $c] = $dataRows[$r][$c]
  }
}

$startCell = $ws4.Cells.Item(1, 1)
$endCell = $ws4.Cells.Item($nRows, $nCols)
$ws4.Range($startCell, $endCell).Value = $arr

# Make the newly added sheet the active tab, matching the saved view state.
$ws4.Activate()
